# edit.ps1
# Applies the "add files for perturbation analysis" change:
#  - Removes the now-unused "baja energia" cell (I4) from the first block
#  - Duplicates the header/body block (rows 2-8) into a new block (rows 12-23)
#    with a new "Presion Solar" / "Cr" / "e" header and 8 extra data rows
#  - Widens column B and re-centers the view on the new content

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "baja energia" label that used to sit in I4 ---
$ws.Range("I4").Value = $null

# --- Widen column B (target OOXML width = 13 characters) ---
$ws.Columns.Item(2).ColumnWidth = 12.1

# --- Build the new second block (rows 12-23), mirroring rows 2-8 ---

# Row 12: top headers (C12 new "Presion Solar" label, D12/H12 stable/unstable)
$ws.Range("C12").Value = "Presion Solar"
$ws.Range("D12").Value = "stable"
$ws.Range("H12").Value = "unstable"

# Row 13: mean dif labels
$ws.Range("D13").Value = "mean dif"
$ws.Range("H13").Value = "mean dif"

# Row 14: x_forward / x_backward labels
$ws.Range("D14").Value = "x_forward"
$ws.Range("H14").Value = "x_backward"

# Row 15: column headers
$ws.Range("C15").Value = "Cr"
$ws.Range("B15").Value = "e"
$ws.Range("D15").Value = "x"
$ws.Range("E15").Value = "y"
$ws.Range("F15").Value = "Vx"
$ws.Range("G15").Value = "Vy"
$ws.Range("H15").Value = "x"
$ws.Range("I15").Value = "y"
$ws.Range("J15").Value = "Vx"
$ws.Range("K15").Value = "Vy"
$ws.Range("L15").Value = "Orbita"

# Data rows 16-23 (parallel arrays: row number, then columns B..K)
# NOTE: values are written in plain decimal (no scientific "E" notation),
# since the host PowerShell parser does not accept exponent literals.
$rowNums = @(16, 17, 18, 19, 20, 21, 22, 23)
$colB = @(0.21, 0.31, 0.41, 0.51, 0.61, 0.71, 0.81, 0.9)
$colC = @(1.21, 1.31, 1.41, 1.51, 1.61, 1.71, 1.81, 1.9)
$colD = @(-0.0001542, -0.000167, -0.0001797, -0.0001925, -0.0002052, -0.0002179, -0.0002307, -0.0002421)
$colE = @(-0.0001795, -0.0001943, -0.0002092, -0.000224, -0.0002388, -0.0002537, -0.0002685, -0.0002819)
$colF = @(-0.0005303, -0.0005745, -0.0006183, -0.0006621, -0.0007059, -0.0007497, -0.0007934, -0.0008328)
$colG = @(-0.0003223, -0.000349, -0.0003756, -0.0004023, -0.000429, -0.0004557, -0.0004823, -0.0005064)
$colH = @(-0.0000741, -0.0000802, -0.0000863, -0.0000924, -0.0000986, -0.0001047, -0.0001108, -0.0001163)
$colI = @(-0.0000562, -0.0000608, -0.0000655, -0.0000701, -0.0000748, -0.0000794, -0.000084, -0.0000882)
$colJ = @(-0.0001759, -0.0001904, -0.0002049, -0.0002195, -0.000234, -0.0002486, -0.0002631, -0.0002762)
$colK = @(-0.000033, -0.0000358, -0.0000385, -0.0000412, -0.000044, -0.0000467, -0.0000494, -0.0000519)

for ($idx = 0; $idx -lt $rowNums.Length; $idx++) {
    $r = $rowNums[$idx]
    $ws.Range("B$r").Value = $colB[$idx]
    $ws.Range("C$r").Value = $colC[$idx]
    $ws.Range("D$r").Value = $colD[$idx]
    $ws.Range("E$r").Value = $colE[$idx]
    $ws.Range("F$r").Value = $colF[$idx]
    $ws.Range("G$r").Value = $colG[$idx]
    $ws.Range("H$r").Value = $colH[$idx]
    $ws.Range("I$r").Value = $colI[$idx]
    $ws.Range("J$r").Value = $colJ[$idx]
    $ws.Range("K$r").Value = $colK[$idx]
    $ws.Range("L$r").Value = "si"

    $bcRange = "B" + $r + ":C" + $r
    $dkRange = "D" + $r + ":K" + $r
    $ws.Range($bcRange).NumberFormat = "General"
    $ws.Range($dkRange).NumberFormat = "0.00E+00"
}

# --- Update the view: scroll down and move the active selection ---
$ws.Range("J25").Select()
$excel.ActiveWindow.ScrollRow = 5
